# Add a question to Q&A file
# New Q&A entry (row 6, previously blank) added to Sheet1: a student question
# (By: Ashkan, Date: 2022-10-02) about training/test image resolutions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Question text (column B) for the new row.
$ws.Range("B6").Value = "Just for clarification: for the resolution of the training image we have: the original image is high resolution (1m) and the label is low resolution (30m). For the test, both images and labels are high resolutions (1m)? I ask because I am sure about the input training (not label) images resolution (I thought It was 30 same as its label but the paper says it is 1m)"

# Asked by (column C).
$ws.Range("C6").Value = "Ashkan"

# Date asked (column D) -- reuse the existing date format from the row above
# so the cell gets the same date-formatted style rather than a brand new one.
$ws.Range("D6").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("D6").Value = 44836

# The row grew taller to fit the wrapped question text.
$ws.Range("A6:H6").RowHeight = 105

# Leave the selection on the newly added row, matching where the edit was made.
$ws.Range("E6").Select()
